$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date strings in column A (rows 3-21): replace "/" separators with "-".
# The leading apostrophe forces Excel to keep the value as literal text instead
# of auto-converting ambiguous day/month strings into a date serial number; the
# Style reset afterwards removes the quote-prefix formatting Excel applies so
# the cell ends up with no explicit style, same as the original.
function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.Formula = "'" + $text
    $rng.Style = "Normal"
}

Set-TextValue "A3"  "28-07-2022"
Set-TextValue "A4"  "01-08-2022"
Set-TextValue "A5"  "04-08-2022"
Set-TextValue "A6"  "08-08-2022"
Set-TextValue "A7"  "11-08-2022"
Set-TextValue "A8"  "15-08-2022"
Set-TextValue "A9"  "18-08-2022"
Set-TextValue "A10" "22-08-2022"
Set-TextValue "A11" "25-08-2022"
Set-TextValue "A12" "29-08-2022"
Set-TextValue "A13" "01-09-2022"
Set-TextValue "A14" "05-09-2022"
Set-TextValue "A15" "08-09-2022"
Set-TextValue "A16" "12-09-2022"
Set-TextValue "A17" "15-09-2022"
Set-TextValue "A18" "19-09-2022"
Set-TextValue "A19" "22-09-2022"
Set-TextValue "A20" "26-09-2022"
Set-TextValue "A21" "29-09-2022"

# Update attendance counters that changed value
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0
